$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) is stored as text in the source sheet (e.g. "62.074.51",
# "409.39"). Plain Value assignment would let Excel auto-coerce simple
# decimal-looking strings into numbers, so NumberFormat is forced to Text
# ("@") immediately before each such write to keep the cell a string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.078.39'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.430.17'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.25'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.94'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.634'
$ws.Range("E7").Value = '  +6.49%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.742'
$ws.Range("E9").Value = '  +7.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +6.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.89'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000229'
$ws.Range("E12").Value = '  +54.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.29'
$ws.Range("E13").Value = '  +10.75%  '
$ws.Range("E15").Value = '  +8.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.972.89'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.398.81'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("E18").Value = '  +7.36%  '
$ws.Range("E19").Value = '  +8.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.081.23'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '454.24'
$ws.Range("E21").Value = '  +45.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.78'
$ws.Range("E22").Value = '  +9.18%  '
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.07'
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("E25").Value = '  +2.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.13'
$ws.Range("E26").Value = '  +11.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.09'
$ws.Range("E27").Value = '  +11.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.76'
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.73'
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.06'
$ws.Range("E31").Value = '  +6.15%  '
$ws.Range("E32").Value = '  -1.10%  '
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.96'
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.17'
$ws.Range("E37").Value = '  +4.92%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("E40").Value = '  +7.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.321'
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("E42").Value = '  -1.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '143.05'
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.26'
$ws.Range("E44").Value = '  +8.94%  '
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.52'
$ws.Range("E46").Value = '  +13.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.69'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.35'
$ws.Range("E48").Value = '  +5.58%  '
$ws.Range("E49").Value = '  +8.96%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.778.63'
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.139'
$ws.Range("E51").Value = '  +15.95%  '
